$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the blank "Summary" row (old row 34); this shifts rows 35-39 up to 34-38
$ws.Rows.Item(34).Delete()

# Re-label column A so each detail row is prefixed with its category (group) name,
# and rename the summary section labels.
$ws.Range("A7").Value  = "     Civilian, New nominations"
$ws.Range("A8").Value  = "     Civilian, Confirmed "
$ws.Range("A9").Value  = "     Civilian, Unconfirmed "
$ws.Range("A10").Value = "     Civilian, Withdrawn "
$ws.Range("A11").Value = "     Civilian, Returned at sine die adjournment "

$ws.Range("A13").Value = "     Civilian (FS, PHS, CG, NOAA), New nominations"
$ws.Range("A14").Value = "     Civilian (FS, PHS, CG, NOAA), Confirmed "
$ws.Range("A15").Value = "     Civilian (FS, PHS, CG, NOAA), Unconfirmed "

$ws.Range("A17").Value = "     Air Force, New nominations"
$ws.Range("A18").Value = "     Air Force, Confirmed"
$ws.Range("A19").Value = "     Air Force, Unconfirmed  "

$ws.Range("A21").Value = "     Army, New nominations"
$ws.Range("A22").Value = "     Army, Confirmed "
$ws.Range("A23").Value = "     Army, Unconfirmed "

$ws.Range("A25").Value = "     Navy, New nominations"
$ws.Range("A26").Value = "     Navy, Confirmed "
$ws.Range("A27").Value = "     Navy, Unconfirmed"
$ws.Range("A28").Value = "     Navy, Returned at sine die adjournment "

$ws.Range("A30").Value = "     Marine Corps, New nominations"
$ws.Range("A31").Value = "     Marine Corps, Confirmed"
$ws.Range("A32").Value = "     Marine Corps, Unconfirmed "
$ws.Range("A33").Value = "     Marine Corps, Withdrawn "

$ws.Range("A34").Value = "Total new nominations"
$ws.Range("A35").Value = "Total confirmed "
$ws.Range("A36").Value = "Total unconfirmed"
$ws.Range("A37").Value = "Total withdrawn "
$ws.Range("A38").Value = "Total returned at sine die adjournment "

# Reset the lingering cell selection (previously on B6, a row that has since
# shifted/changed) back to the sheet's default top-left cell.
$ws.Range("A1").Select()
